# quarterly_seprated.xlsx -- "add monte_carlo and update database"
#
# The commit updates the "فصل سوم منتهی به 1400/09" quarter column (column J)
# of the Overview sheet: some quarterly sales-volume figures are blanked out
# (replaced with the existing "-" placeholder string already used elsewhere
# in the sheet for missing data), and several quarterly amount/cost/profit
# figures are revised to new computed totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Cells whose value becomes the "-" placeholder (no data for the quarter) ---
$dashCells = @("J11","J12","J13","J14","J18","J19","J20","J26","J34","J35","J36","J37","J41","J42","J43","J49")
foreach ($cellRef in $dashCells) {
    $ws.Range($cellRef).Value = "-"
}

# --- Cells whose numeric value is reset to 0 ---
$zeroCells = @("J16","J22","J27","J39","J45","J50")
foreach ($cellRef in $zeroCells) {
    $ws.Range($cellRef).Value = 0
}

# --- Cells whose numeric value is revised to a new figure ---
$ws.Range("J59").Value  = 84588426
$ws.Range("J73").Value  = -2424058
$ws.Range("J74").Value  = -16633831
$ws.Range("J75").Value  = -1303499
$ws.Range("J76").Value  = -18875935
$ws.Range("J78").Value  = -39237323
$ws.Range("J80").Value  = -21758346
$ws.Range("J82").Value  = -19008203
$ws.Range("J84").Value  = -41635738
$ws.Range("J89").Value  = -80873061
$ws.Range("J96").Value  = 499837
$ws.Range("J97").Value  = 12187837
$ws.Range("J98").Value  = 858412
$ws.Range("J99").Value  = 8062719
$ws.Range("J101").Value = 21608805
$ws.Range("J103").Value = 15120917
$ws.Range("J105").Value = 10563409
$ws.Range("J107").Value = 26224849
$ws.Range("J110").Value = 47833654
